$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'youth black knee pads'
    'basketball clothes for boys'
    'men capri'
    'black baseball pants adult small'
    'under pants for men'
    'by design knee pads'
    'baseball knee high pants mens'
    'boys youth compression tights'
    'mens small running tights'
    'youth large softball pants'
    'small basketballs'
    'cycling capri'
    'mens black basketball shorts'
    'knee protector for running'
    'knees protection'
    'dry fit baseball pants'
    'girls black baseball pants'
    'mens baseball sliding shorts'
    'wrestling equipment'
    'the knee pads'
    'football compression'
    'soccer guards for men'
    'spandex leggings for boys'
    'basketball pants men'
    'hunting knee pads'
    'combat pants knee pads'
    'nike compression basketball pants'
    'seamless capri leggings'
    'basketball knee pads adidas'
    'women black leggings'
    'red basketball knee pads'
    'underarmour mens leggings'
    'red nike compression pants men'
    'men compression pants long'
    'knee tights'
    'knee protector basketball'
    'youth football pants with pads'
    'little boys compression tights'
    'wrestling tights for boys'
    'youth boys baseball pants'
    'padded basketball tights'
    'basketball clothes youth'
    'leggings boys'
    'youth compression pants'
    'knee pad youth'
    'best knee pads for volleyball'
    'athletic pants for men big and tall'
    'leg compression basketball'
    'compression shorts baseball'
    'pants for men sports'
    'knee pads large'
    'gym tights for men'
    'soccer gear for men'
    'large tall athletic pants men'
    'youth football pads'
    'padded compression shorts'
    'youth girls softball pants black'
    'athletic pants for men'
    'patella knee'
    'youth softball pants'
    'basketball knee sleeves'
    'knee pads for crossfit'
    'mens nike dri fit compression pants'
    'volleyball youth knee pads'
    'troll knee pads'
    'elastic knee pads'
    'men leggings adidas'
    'snowmobile knee pads'
    'knee pads military'
    'knee pads sleeve'
    'knee pads mma'
    'under armour compression pants men'
    'mcdavid knee pads basketball'
    'black leggings xsmall'
    'tesla compression pants'
    'fox knee pads'
    'biking knee pads'
    'bmx knee pads'
    'mcdavid knee pad'
    'pant with knee pads'
    'green knee pads'
    'knee pads tights'
    'warm compression pants'
    'dancer knee pads'
    'navy compression pants'
    'mizuno knee pad'
    'mens pants with knee pads'
    'black athletic capri'
    'tactical knee pad'
    'men compression pants blue'
    'leggins for men sport'
    'knee pads ski'
    'knee pad snowboard'
    'knee pads mcdavid basketball'
    'pantalones con rodilleras'
    'pantalon con rodilleras'
    'cat knee pad pants'
    'knee padded pants men'
    'knee pad pants men'
    'caterpillar knee pad pants'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
